{"js": "// Update the title date and the 25 division-problem answers in the single\n// table (5 populated rows x 5 columns, at table-row indices 0,4,8,12,16;\n// the intervening rows are blank spacer rows).\n\n// 1) Title paragraph: \"2023-11-25 Saturday\" -> \"2023-11-26 Sunday\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\nif (titlePara.text.trim() === \"2023-11-25 Saturday\") {\n  // Replace via the paragraph's first run range text to keep formatting.\n  const range = titlePara.getRange();\n  range.insertText(\"2023-11-26 Sunday\", Word.InsertLocation.replace);\n}\n\n// 2) Table cell values (row index, col index) -> new text.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst updates = [\n  [0, 0, \"18\u00f77=2, 4\"],\n  [0, 1, \"11\u00f77=1, 4\"],\n  [0, 2, \"26\u00f76=4, 2\"],\n  [0, 3, \"12\u00f75=2, 2\"],\n  [0, 4, \"37\u00f73=12, 1\"],\n  [4, 0, \"89\u00f76=14, 5\"],\n  [4, 1, \"98\u00f74=24, 2\"],\n  [4, 2, \"93\u00f77=13, 2\"],\n  [4, 3, \"31\u00f78=3, 7\"],\n  [4, 4, \"64\u00f78=8, 0\"],\n  [8, 0, \"60\u00f73=20, 0\"],\n  [8, 1, \"78\u00f74=19, 2\"],\n  [8, 2, \"54\u00f77=7, 5\"],\n  [8, 3, \"87\u00f75=17, 2\"],\n  [8, 4, \"54\u00f79=6, 0\"],\n  [12, 0, \"15\u00f77=2, 1\"],\n  [12, 1, \"70\u00f75=14, 0\"],\n  [12, 2, \"32\u00f78=4, 0\"],\n  [12, 3, \"15\u00f75=3, 0\"],\n  [12, 4, \"82\u00f79=9, 1\"],\n  [16, 0, \"89\u00f78=11, 1\"],\n  [16, 1, \"55\u00f73=18, 1\"],\n  [16, 2, \"79\u00f78=9, 7\"],\n  [16, 3, \"44\u00f78=5, 4\"],\n  [16, 4, \"48\u00f78=6, 0\"],\n];\n\nfor (const [row, col, text] of updates) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the title date and the 25 division-problem answers in the single\n# table (5 populated rows x 5 columns, at 1-based table row indices\n# 1,5,9,13,17; the intervening rows are blank spacer rows).\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2023-11-25 Saturday\" -> \"2023-11-26 Sunday\"\n$titlePara = $d.Paragraphs.First\nif ($titlePara.Range.Text.Trim() -eq \"2023-11-25 Saturday\") {\n    $titlePara.Range.Text = \"2023-11-26 Sunday\"\n}\n\n# 2) Table cell values (row, column) -> new text (1-based, as COM expects).\n$tbl = $d.Tables.Item(1)\n\n$updates = @(\n    @(1, 1, \"18\u00f77=2, 4\"),\n    @(1, 2, \"11\u00f77=1, 4\"),\n    @(1, 3, \"26\u00f76=4, 2\"),\n    @(1, 4, \"12\u00f75=2, 2\"),\n    @(1, 5, \"37\u00f73=12, 1\"),\n    @(5, 1, \"89\u00f76=14, 5\"),\n    @(5, 2, \"98\u00f74=24, 2\"),\n    @(5, 3, \"93\u00f77=13, 2\"),\n    @(5, 4, \"31\u00f78=3, 7\"),\n    @(5, 5, \"64\u00f78=8, 0\"),\n    @(9, 1, \"60\u00f73=20, 0\"),\n    @(9, 2, \"78\u00f74=19, 2\"),\n    @(9, 3, \"54\u00f77=7, 5\"),\n    @(9, 4, \"87\u00f75=17, 2\"),\n    @(9, 5, \"54\u00f79=6, 0\"),\n    @(13, 1, \"15\u00f77=2, 1\"),\n    @(13, 2, \"70\u00f75=14, 0\"),\n    @(13, 3, \"32\u00f78=4, 0\"),\n    @(13, 4, \"15\u00f75=3, 0\"),\n    @(13, 5, \"82\u00f79=9, 1\"),\n    @(17, 1, \"89\u00f78=11, 1\"),\n    @(17, 2, \"55\u00f73=18, 1\"),\n    @(17, 3, \"79\u00f78=9, 7\"),\n    @(17, 4, \"44\u00f78=5, 4\"),\n    @(17, 5, \"48\u00f78=6, 0\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $cell = $tbl.Cell($row, $col)\n    $cell.Range.Text = $text\n}\n"}
